# Update countries & provincias Spain
# Applies the 12-Sep-2020 18:11 data refresh: new totals for several countries,
# which also shifts a few countries rank (so the same spreadsheet row now
# shows a different country name), plus the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" banner in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 18:11"

# Rank-shuffle: same row, new country label (Casos totales column still sorted desc)
$countryRenames = [ordered]@{
    16  = "Reino Unido"
    17  = "Francia"
    22  = "Irak"
    23  = "Italia"
    120 = "Uganda"
    121 = "Cabo Verde"
    122 = "Cuba"
    123 = "Ruanda"
    124 = "Surinam"
    143 = "Birmania"
    144 = "Sudan del Sur"
}
foreach ($row in $countryRenames.Keys) {
    $ws.Range("A$row").Value = $countryRenames[$row]
}

# Updated per-country stats: B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
$statUpdates = @{
    4   = @{ B=6647808; C=11561; D=3919483; E=2530661; F=0; G=243; H=197664 }
    14  = @{ B=432666;  C=2131;  D=404919;  E=15852;   F=0; G=45;  H=11895 }
    16  = @{ B=365174;  C=3497;  D=0;       E=0;       F=0; G=9;   H=41623 }
    17  = @{ B=363350;  C=0;     D=89059;   E=243398;  F=0; G=0;   H=30893 }
    22  = @{ B=286778;  C=4106;  D=221283;  E=57554;   F=0; G=60;  H=7941  }
    23  = @{ B=286297;  C=1501;  D=213191;  E=37503;   F=0; G=6;   H=35603 }
    55  = @{ D=56699; E=631 }
    120 = @{ B=4703; C=326; D=1998; E=2653; F=0; G=3; H=52  }
    121 = @{ B=4651; C=0;   D=4076; E=531;  F=0; G=0; H=44  }
    122 = @{ B=4593; C=0;   D=3844; E=643;  F=0; G=0; H=106 }
    123 = @{ B=4534; C=0;   D=2450; E=2062; F=0; G=0; H=22  }
    124 = @{ B=4529; C=0;   D=3747; E=689;  F=0; G=0; H=93  }
    136 = @{ B=3062; C=117; D=2156; E=884;  F=0; G=1; H=22  }
    143 = @{ B=2595; C=173; D=676;  E=1903; F=0; G=2; H=16  }
    144 = @{ B=2568; C=0;   D=1290; E=1229; F=0; G=0; H=49  }
    193 = @{ B=139;  C=1;   D=1;    E=3;    F=0; G=0; H=0   }
}
foreach ($row in $statUpdates.Keys) {
    $rowVals = $statUpdates[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
